$wb = $excel.ActiveWorkbook

# "展览" (Exhibition) sheet
$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F3").Value = 396
$wsExhibition.Range("F4").Value = 5005
$wsExhibition.Range("F5").Value = 37
$wsExhibition.Range("F8").Value = 496

# "全部类型" (All Types) sheet
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 396
$wsAll.Range("F4").Value = 5005
$wsAll.Range("F6").Value = 37
$wsAll.Range("F10").Value = 496
